$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy header style (style index 1) from AC1 to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

# Fill in win/loss/tie record for every data row (2 through 58)
$lastRow = 58
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 76  # AD
    $ws.Cells.Item($r, 31).Value = 86  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
